# Update metricas_recorrencia_trimestral row 21 (2025Q3) with refreshed data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = 258
$ws.Range("D21").Value = 228
$ws.Range("E21").Value = 30
$ws.Range("F21").Value = 65.32951289398281
